$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "Teachers") to match the export file name.
$ws.Name = "giaovien.xlsx"

# Fix header text: "School ID" -> "SchoolID" (also drops the now-unused
# empty shared-string entry automatically since it is no longer referenced).
$ws.Range("C1").Value = "SchoolID"

# Column C best-fit width shrinks now that the header text is one
# character shorter; nudge the stored column width down to match.
$ws.Columns.Item(3).ColumnWidth = 10.67
